$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 1411
$wsExpo.Range("F5").Value = 673

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1411
$wsAll.Range("F6").Value = 673
